$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 833.3333
$ws.Range("I18").Value = 833.3333
$ws.Range("K18").Value = 833.3333
$ws.Range("M18").Value = -549.3333
$ws.Range("H74").Value = 6843.7334
$ws.Range("I74").Value = 5923.364
$ws.Range("K74").Value = 5923.364
$ws.Range("M74").Value = -4987.364
$ws.Range("H77").Value = 6843.7334
$ws.Range("I77").Value = 5923.364
$ws.Range("K77").Value = 29616.82
$ws.Range("M77").Value = -24936.82
$ws.Range("H100").Value = 6825.1035
$ws.Range("I100").Value = 1306.8
$ws.Range("K100").Value = 1306.8
$ws.Range("M100").Value = -765.8
$ws.Range("H103").Value = 438.9
$ws.Range("I103").Value = 231.88889
$ws.Range("J103").Value = 608.2727
$ws.Range("K103").Value = 695.6666700000001
$ws.Range("L103").Value = 1824.8181
$ws.Range("M103").Value = -109.6666700000001
$ws.Range("N103").Value = -2996.8181
$ws.Range("H137").Value = 2337.724
$ws.Range("I137").Value = 2375.6785
$ws.Range("J137").Value = 1275
$ws.Range("K137").Value = 7127.0355
$ws.Range("L137").Value = 3825
$ws.Range("M137").Value = -4577.0355
$ws.Range("N137").Value = -8925

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1368.1904
$ws.Range("I2").Value = 1555.4286
$ws.Range("K2").Value = 1555.4286
$ws.Range("M2").Value = -1442.4286
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()
$ws.Range("H43").Value = 13636.125
$ws.Range("I43").Value = 10171
$ws.Range("J43").Value = 14791.167
$ws.Range("K43").Value = 10171
$ws.Range("L43").Value = 14791.167
$ws.Range("M43").Value = -9858
$ws.Range("N43").Value = -15417.167
$ws.Range("H116").Value = 1368.1904
$ws.Range("I116").Value = 1555.4286
$ws.Range("K116").Value = 1555.4286
$ws.Range("M116").Value = 738.5714
$ws.Range("H132").Value = 1797.8695
$ws.Range("I132").Value = 1648.9756
$ws.Range("K132").Value = 4946.9268
$ws.Range("M132").Value = -2416.9268

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1368.1904
$ws.Range("I3").Value = 1555.4286
$ws.Range("K3").Value = 1555.4286
$ws.Range("M3").Value = -1441.4286
$ws.Range("H20").Value = 1397.9111
$ws.Range("I20").Value = 1292.8518
$ws.Range("J20").Value = 1555.5
$ws.Range("K20").Value = 1292.8518
$ws.Range("L20").Value = 1555.5
$ws.Range("M20").Value = -1045.8518
$ws.Range("N20").Value = -2049.5
$ws.Range("H109").Value = 99638.5
$ws.Range("J109").Value = 99638.5
$ws.Range("L109").Value = 99638.5
$ws.Range("N109").Value = -102412.5
$ws.Range("H115").Value = 63999
$ws.Range("J115").Value = 63999
$ws.Range("L115").Value = 63999
$ws.Range("N115").Value = -67133
$ws.Range("H123").Value = 45000
$ws.Range("J123").Value = 45000
$ws.Range("L123").Value = 45000
$ws.Range("N123").Value = -54800

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1779.8572
$ws.Range("I16").Value = 1409.8334
$ws.Range("K16").Value = 1409.8334
$ws.Range("M16").Value = -1122.8334
$ws.Range("H94").Value = 1545.9445
$ws.Range("J94").Value = 1718.5
$ws.Range("L94").Value = 1718.5
$ws.Range("N94").Value = -2620.5
$ws.Range("H112").Value = 100000.5
$ws.Range("J112").Value = 100000.5
$ws.Range("L112").Value = 100000.5
$ws.Range("N112").Value = -102954.5
$ws.Range("H113").Value = 1779.8572
$ws.Range("I113").Value = 1409.8334
$ws.Range("K113").Value = 1409.8334
$ws.Range("M113").Value = 760.1666
$ws.Range("H141").Value = 183347.89
$ws.Range("J141").Value = 183347.89
$ws.Range("L141").Value = 183347.89
$ws.Range("N141").Value = -193707.89

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 2283
$ws.Range("I34").Value = 200
$ws.Range("J34").Value = 3324.5
$ws.Range("K34").Value = 600
$ws.Range("L34").Value = 9973.5
$ws.Range("M34").Value = -516
$ws.Range("N34").Value = -10141.5
$ws.Range("H39").Value = 6306.857
$ws.Range("I39").Value = 399.33334
$ws.Range("J39").Value = 10737.5
$ws.Range("K39").Value = 1198.00002
$ws.Range("L39").Value = 32212.5
$ws.Range("M39").Value = -904.0000199999999
$ws.Range("N39").Value = -32800.5
$ws.Range("H57").Value = 5000
$ws.Range("I57").Value = 5000
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 15000
$ws.Range("L57").Value = 0
$ws.Range("M57").Value = -14441
$ws.Range("N57").ClearContents()
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("M76").ClearContents()
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("M79").ClearContents()
$ws.Range("N79").ClearContents()
$ws.Range("H132").Value = 1254
$ws.Range("I132").Value = 997.3333
$ws.Range("J132").Value = 1474
$ws.Range("K132").Value = 8975.9997
$ws.Range("L132").Value = 13266
$ws.Range("M132").Value = -6445.9997
$ws.Range("N132").Value = -18326
$ws.Range("H134").Value = 3493.7778
$ws.Range("I134").Value = 2492
$ws.Range("J134").Value = 7000
$ws.Range("K134").Value = 7476
$ws.Range("L134").Value = 21000
$ws.Range("M134").Value = -2406
$ws.Range("N134").Value = -31140
$ws.Range("H140").Value = 35286.383
$ws.Range("I140").Value = 46021.76
$ws.Range("J140").Value = 5465.8887
$ws.Range("K140").Value = 138065.28
$ws.Range("L140").Value = 16397.6661
$ws.Range("M140").Value = -132885.28
$ws.Range("N140").Value = -26757.6661

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5047.952
$ws.Range("I70").Value = 4114.3335
$ws.Range("J70").Value = 5748.1665
$ws.Range("K70").Value = 4114.3335
$ws.Range("L70").Value = 5748.1665
$ws.Range("M70").Value = -3844.3335
$ws.Range("N70").Value = -6288.1665
$ws.Range("H73").Value = 5047.952
$ws.Range("I73").Value = 4114.3335
$ws.Range("J73").Value = 5748.1665
$ws.Range("K73").Value = 4114.3335
$ws.Range("L73").Value = 5748.1665
$ws.Range("M73").Value = -3178.3335
$ws.Range("N73").Value = -7620.1665
$ws.Range("H132").Value = 11115646
$ws.Range("I132").Value = 18522852
$ws.Range("J132").Value = 4836.4165
$ws.Range("K132").Value = 55568556
$ws.Range("L132").Value = 14509.2495
$ws.Range("M132").Value = -55566026
$ws.Range("N132").Value = -19569.2495
$ws.Range("H136").Value = 50058.75
$ws.Range("J136").Value = 50058.75
$ws.Range("L136").Value = 150176.25
$ws.Range("N136").Value = -155276.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4097.36
$ws.Range("I40").Value = 3922.8635
$ws.Range("K40").Value = 3922.8635
$ws.Range("M40").Value = -3786.8635
$ws.Range("H133").Value = 70127.75
$ws.Range("I133").Value = 67296
$ws.Range("J133").Value = 71071.664
$ws.Range("K133").Value = 67296
$ws.Range("L133").Value = 71071.664
$ws.Range("M133").Value = -64766
$ws.Range("N133").Value = -76131.664

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 83348610
$ws.Range("I2").Value = 17277.5
$ws.Range("J2").Value = 166679940
$ws.Range("K2").Value = 17277.5
$ws.Range("L2").Value = 166679940
$ws.Range("M2").Value = -17165.5
$ws.Range("N2").Value = -166680164
$ws.Range("H4").Value = 90917050
$ws.Range("I4").Value = 4083.3333
$ws.Range("J4").Value = 200012590
$ws.Range("K4").Value = 4083.3333
$ws.Range("L4").Value = 200012590
$ws.Range("M4").Value = -3970.3333
$ws.Range("N4").Value = -200012816
$ws.Range("H107").Value = 722.3200000000001
$ws.Range("I107").Value = 416
$ws.Range("J107").Value = 1112.1818
$ws.Range("K107").Value = 1248
$ws.Range("L107").Value = 3336.5454
$ws.Range("M107").Value = 672
$ws.Range("N107").Value = -7176.5454
